$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.130.66'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.655.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.68%  '
$ws.Range("E4").Value = '  -0.47%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5260'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.39%  '
$ws.Range("E7").Value = '  -0.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2686'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06377'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.58'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07690'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.615'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.578.57'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.883.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5638'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8261'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.71'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.115.98'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.67%  '
$ws.Range("E19").Value = '  -0.46%  '
$ws.Range("E20").Value = '  -0.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '190.41'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.003'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.97%  '
$ws.Range("E24").Value = '  -0.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1203'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.271'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.01'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.97%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.522'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05655'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.26%  '
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("E32").Value = '  -0.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.387'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.583'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.799'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9491'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.409'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5782'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("E39").Value = '  -0.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.976'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8371'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.023.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.76%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.42'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.793.63'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '58.43'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₈105'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.004'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05340'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.052'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4341'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.60%  '
